$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1544.4445
$ws.Range("I69").Value = 2300
$ws.Range("K69").Value = 6900
$ws.Range("M69").Value = -6026

$ws.Range("H72").Value = 1544.4445
$ws.Range("I72").Value = 2300
$ws.Range("K72").Value = 20700
$ws.Range("M72").Value = -16332

$ws.Range("H86").Value = 7343.1177
$ws.Range("I86").Value = 1535.3
$ws.Range("K86").Value = 1535.3
$ws.Range("M86").Value = -412.3

$ws.Range("H89").Value = 7343.1177
$ws.Range("I89").Value = 1535.3
$ws.Range("K89").Value = 7676.5
$ws.Range("M89").Value = -2060.5

$ws.Range("H107").Value = 1438.2354
$ws.Range("I107").Value = 1573.0769
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1573.0769
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 346.9231
$ws.Range("N107").Value = -4840

$ws.Range("H129").Value = 855
$ws.Range("J129").Value = 860
$ws.Range("L129").Value = 2580
$ws.Range("N129").Value = -12580

$ws.Range("H132").Value = 3541.08
$ws.Range("I132").Value = 3546.524
$ws.Range("J132").Value = 3512.5
$ws.Range("K132").Value = 10639.572
$ws.Range("L132").Value = 10537.5
$ws.Range("M132").Value = -8109.572
$ws.Range("N132").Value = -15597.5

$ws.Range("H135").Value = 22733118
$ws.Range("I135").Value = 689
$ws.Range("J135").Value = 125029050
$ws.Range("K135").Value = 6201
$ws.Range("L135").Value = 1125261450
$ws.Range("M135").Value = -3666
$ws.Range("N135").Value = -1125266520

$ws.Range("H137").Value = 85468.75
$ws.Range("I137").Value = 2025.5
$ws.Range("K137").Value = 6076.5
$ws.Range("M137").Value = -3526.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2319
$ws.Range("I2").Value = 1495.4615
$ws.Range("K2").Value = 1495.4615
$ws.Range("M2").Value = -1382.4615

$ws.Range("H74").Value = 2156.5186
$ws.Range("I74").Value = 2104.8462
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 2104.8462
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -1230.8462
$ws.Range("N74").Value = -5248

$ws.Range("H77").Value = 2156.5186
$ws.Range("I77").Value = 2104.8462
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 10524.231
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -6156.231
$ws.Range("N77").Value = -26236

$ws.Range("H97").Value = 2226.5386
$ws.Range("I97").Value = 1826.5
$ws.Range("K97").Value = 1826.5
$ws.Range("M97").Value = -1330.5

$ws.Range("H116").Value = 2319
$ws.Range("I116").Value = 1495.4615
$ws.Range("K116").Value = 1495.4615
$ws.Range("M116").Value = 798.5385000000001

$ws.Range("H122").Value = 2348.3333
$ws.Range("I122").Value = 2518
$ws.Range("K122").Value = 7554
$ws.Range("M122").Value = -5104

$ws.Range("H138").Value = 31143.2
$ws.Range("J138").Value = 31143.2
$ws.Range("L138").Value = 31143.2
$ws.Range("N138").Value = -41423.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2319
$ws.Range("I3").Value = 1495.4615
$ws.Range("K3").Value = 1495.4615
$ws.Range("M3").Value = -1381.4615

$ws.Range("H100").Value = 26814.334
$ws.Range("J100").Value = 26814.334
$ws.Range("L100").Value = 26814.334
$ws.Range("N100").Value = -28978.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2701.4644
$ws.Range("I31").Value = 2266.842
$ws.Range("J31").Value = 3619
$ws.Range("K31").Value = 2266.842
$ws.Range("L31").Value = 3619
$ws.Range("M31").Value = -1971.842
$ws.Range("N31").Value = -4209

$ws.Range("H34").Value = 2701.4644
$ws.Range("I34").Value = 2266.842
$ws.Range("J34").Value = 3619
$ws.Range("K34").Value = 2266.842
$ws.Range("L34").Value = 3619
$ws.Range("M34").Value = -2064.842
$ws.Range("N34").Value = -4023

$ws.Range("H96").Value = 4262
$ws.Range("J96").Value = 4262
$ws.Range("L96").Value = 4262
$ws.Range("N96").Value = -9754

$ws.Range("H122").Value = 911.9
$ws.Range("I122").Value = 911.9
$ws.Range("K122").Value = 2735.7
$ws.Range("M122").Value = -285.6999999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 4667.8
$ws.Range("I63").Value = 1905.5
$ws.Range("J63").Value = 6509.3335
$ws.Range("K63").Value = 5716.5
$ws.Range("L63").Value = 19528.0005
$ws.Range("M63").Value = -4967.5
$ws.Range("N63").Value = -21026.0005

$ws.Range("H66").Value = 4667.8
$ws.Range("I66").Value = 1905.5
$ws.Range("J66").Value = 6509.3335
$ws.Range("K66").Value = 17149.5
$ws.Range("L66").Value = 58584.0015
$ws.Range("M66").Value = -13405.5
$ws.Range("N66").Value = -66072.0015

$ws.Range("H92").Value = 866.55554
$ws.Range("J92").Value = 1866.3334
$ws.Range("L92").Value = 5599.0002
$ws.Range("N92").Value = -8095.0002

$ws.Range("H131").Value = 756.27
$ws.Range("J131").Value = 777.4421
$ws.Range("L131").Value = 2332.3263
$ws.Range("N131").Value = -12412.3263

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H102").Value = 1558.1945
$ws.Range("I102").Value = 1636.25
$ws.Range("J102").Value = 1285
$ws.Range("K102").Value = 1636.25
$ws.Range("L102").Value = 1285
$ws.Range("M102").Value = -14.25
$ws.Range("N102").Value = -4529

$ws.Range("H122").Value = 2350.625
$ws.Range("I122").Value = 2359.8
$ws.Range("J122").Value = 2335.3333
$ws.Range("K122").Value = 7079.400000000001
$ws.Range("L122").Value = 7005.999899999999
$ws.Range("M122").Value = -4629.400000000001
$ws.Range("N122").Value = -11905.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1512126
$ws.Range("I122").Value = 3924398.8
$ws.Range("K122").Value = 11773196.4
$ws.Range("M122").Value = -11770746.4

$ws.Range("H141").Value = 67000
$ws.Range("J141").Value = 67000
$ws.Range("L141").Value = 67000
$ws.Range("N141").Value = -77360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2208.1333
$ws.Range("I122").Value = 2019.5834
$ws.Range("J122").Value = 2962.3333
$ws.Range("K122").Value = 6058.7502
$ws.Range("L122").Value = 8886.999899999999
$ws.Range("M122").Value = -3608.7502
$ws.Range("N122").Value = -13786.9999

Write-Output "edits applied"
